$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.718.35"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "3.306.78"
$ws.Range("E3").Value = "  +1.88%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.Value = "'605.38"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.84%  "
$c = $ws.Range("D6")
$c.Value = "'141.41"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.304.21"
$ws.Range("E8").Value = "  +1.99%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +1.66%  "
$c = $ws.Range("D11")
$c.Value = "'5.55"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.61%  "
$c = $ws.Range("D12")
$c.Value = "'0.469"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  +0.38%  "
$c = $ws.Range("D14")
$c.Value = "'34.92"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "3.850.19"
$ws.Range("E15").Value = "  +2.00%  "
$ws.Range("E16").Value = "  +0.62%  "
$ws.Range("D17").Value = "3.305.86"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "63.805.07"
$ws.Range("E18").Value = "  +0.70%  "
$c = $ws.Range("D19")
$c.Value = "'6.87"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "
$c = $ws.Range("D20")
$c.Value = "'480.12"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.26%  "
$c = $ws.Range("D21")
$c.Value = "'14.03"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.94%  "
$c = $ws.Range("D22")
$c.Value = "'0.740"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.09%  "
$c = $ws.Range("D23")
$c.Value = "'7.99"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.53%  "
$c = $ws.Range("D24")
$c.Value = "'13.94"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +5.73%  "
$c = $ws.Range("D25")
$c.Value = "'85.25"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("E28").Value = "  +0.06%  "
$c = $ws.Range("D29")
$c.Value = "'7.26"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -1.29%  "
$c = $ws.Range("D30")
$c.Value = "'8.18"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +1.08%  "
$c = $ws.Range("D32")
$c.Value = "'28.77"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +4.47%  "
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("E35").Value = "  +0.98%  "
$c = $ws.Range("D36")
$c.Value = "'6.05"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D37")
$c.Value = "'52.41"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0747"
$ws.Range("E38").Value = "  +5.12%  "
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "3.111.31"
$ws.Range("E40").Value = "  +4.62%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D41")
$c.Value = "'430.46"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  +7.82%  "
$c = $ws.Range("D43")
$c.Value = "'2.74"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("E44").Value = "  -0.71%  "
$c = $ws.Range("D45")
$c.Value = "'0.264"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("E46").Value = "  +2.63%  "
$c = $ws.Range("D47")
$c.Value = "'36.82"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +9.13%  "
$c = $ws.Range("D48")
$c.Value = "'26.40"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.99%  "
$c = $ws.Range("D50")
$c.Value = "'126.24"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("E51").Value = "  -0.22%  "
